# Natmi following Dr Hou advice
# Update ligand/receptor expression statistics (and dependent edge-weight
# columns) for rows 2-7 now that the ligand/receptor-expressing cell counts
# changed from 1 to 3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = @{
        E = 3; G = 13.79659733333333; H = 41.389792
        I = 0.9485830781324925; J = 0.9485830781324925
        K = 3; M = 4.675378666666666; N = 14.026136
        O = 0.03681964474327726; P = 0.03681964474327726
        Q = 64.50431684485687; R = 580.5388516037119
        S = 0.03492649194632279; T = 0.03492649194632279
    }
    3  = @{
        E = 3; G = 13.79659733333333; H = 41.389792
        I = 0.9485830781324925; J = 0.9485830781324925
        K = 3; M = 82.95722966666666; N = 248.871689
        O = 0.653306596744776; P = 0.653306596744776
        Q = 1144.527493599854; R = 10300.74744239869
        S = 0.6197155825044226; T = 0.6197155825044226
    }
    4  = @{
        E = 3; G = 13.79659733333333; H = 41.389792
        I = 0.9485830781324925; J = 0.9485830781324925
        K = 3; M = 39.34793966666667; N = 118.043819
        O = 0.3098737585119468; P = 0.3098737585119468
        Q = 542.8676794772942; R = 4885.809115295648
        S = 0.2939410036817472; T = 0.2939410036817471
    }
    5  = @{
        E = 3; G = 0.7478296666666666; H = 2.243489
        I = 0.05141692186750751; J = 0.05141692186750751
        K = 3; M = 4.675378666666666; N = 14.026136
        O = 0.03681964474327726; P = 0.03681964474327726
        Q = 3.496386869833777; R = 31.467481828504
        S = 0.001893152796954471; T = 0.001893152796954471
    }
    6  = @{
        E = 3; G = 0.7478296666666666; H = 2.243489
        I = 0.05141692186750751; J = 0.05141692186750751
        K = 3; M = 82.95722966666666; N = 248.871689
        O = 0.653306596744776; P = 0.653306596744776
        Q = 62.03787740921344; R = 558.340896682921
        S = 0.03359101424035339; T = 0.03359101424035339
    }
    7  = @{
        E = 3; G = 0.7478296666666666; H = 2.243489
        I = 0.05141692186750751; J = 0.05141692186750751
        K = 3; M = 39.34793966666667; N = 118.043819
        O = 0.3098737585119468; P = 0.3098737585119468
        Q = 29.42555660494344; R = 264.830009444491
        S = 0.01593275483019966; T = 0.01593275483019966
    }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
